$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 ("tutorial wave" row), shifting
# everything else down by one, and fill in the new "tutorial" entry.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "Tutorial text"
$ws.Range("B2").Value = "tutorial"
$ws.Range("C2").Value = "Tutorial"
$ws.Range("D2").Value = "Tutoriel"
$ws.Range("E2").Value = "Tutorial"
$ws.Range("F2").Value = "チュートリアル。"
$ws.Range("G2").Value = "教程"

# Match style used on sibling data rows
$ws.Range("A3:G3").Copy()
$ws.Range("A2:G2").PasteSpecial(-4122)

# Update the active selection to match the new edit location.
$ws.Range("C3").Select()
